$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($sheet, $row1, $row2, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $v1 = $sheet.Cells.Item($row1, $c).Value()
        $v2 = $sheet.Cells.Item($row2, $c).Value()
        $sheet.Cells.Item($row1, $c).Value = $v2
        $sheet.Cells.Item($row2, $c).Value = $v1
    }
}

# Rows 58 and 59 had their match data (columns F..V) swapped
Swap-RowRange $ws 58 59 6 22

# Rows 64 and 65 had their match data (columns F..V) swapped
Swap-RowRange $ws 64 65 6 22

# Append new row 140 (Indice 139) - Mallorca vs Cadiz CF
$newRow = 140

# Copy formatting for the index column (A) and the date column (E) from the
# last existing data row (139) so the new row matches the sheet's style.
$ws.Cells.Item(139, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item(139, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 139
$ws.Cells.Item($newRow, 2).Value = "spain"
$ws.Cells.Item($newRow, 3).Value = "laliga"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45259.875
$ws.Cells.Item($newRow, 6).Value = "Mallorca"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Cadiz CF"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 1.87
$ws.Cells.Item($newRow, 11).Value = "29/10/2023 11:02"
$ws.Cells.Item($newRow, 12).Value = 1.87
$ws.Cells.Item($newRow, 13).Value = "29/11/2023 20:58"
$ws.Cells.Item($newRow, 14).Value = 3.29
$ws.Cells.Item($newRow, 15).Value = "29/10/2023 11:02"
$ws.Cells.Item($newRow, 16).Value = 3.33
$ws.Cells.Item($newRow, 17).Value = "29/11/2023 20:58"
$ws.Cells.Item($newRow, 18).Value = 4.96
$ws.Cells.Item($newRow, 19).Value = "29/10/2023 11:02"
$ws.Cells.Item($newRow, 20).Value = 5.15
$ws.Cells.Item($newRow, 21).Value = "29/11/2023 20:58"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/spain/laliga/mallorca-cadiz/t4abRzTb/"

Write-Output "edit complete"
